$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D11: "asd" -> "= "
$ws.Range("D11").Value = "= "

# Row 54 used to hold: A54="a", B54="= ", C54="plus(2*a+1,a+25*3)"
# It becomes a single combined formula-like string in A54, and B54/C54 are cleared.
$ws.Range("A54").Value = "a = plus(2*a+1,a+25*3)"
$ws.Range("B54").ClearContents()
$ws.Range("C54").ClearContents()

# New rows with additional sample expressions.
$ws.Range("A55").Value = "a--100"
$ws.Range("B55").Value = ""
$ws.Range("A56").Value = "b--100"
$ws.Range("A57").Value = "a++100"

# Move the selection/active cell to A58, matching the post-edit cursor position.
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A58").Select()
